# Edit script: applies data-correction swaps/rotations to the Artfynd sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51   # column AY

function Get-RowValues($row) {
    $vals = New-Object object[] $lastCol
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals[$c - 1] = $ws.Cells.Item($row, $c).Value
    }
    return ,$vals
}

function Set-RowValues($row, $vals) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}

function Swap-Rows($r1, $r2) {
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# 1. Simple value correction
$ws.Range("B2").Value = 91833

# 2. Full swap of row 5 and row 7
Swap-Rows 5 7

# 3. Full swap of row 12 and row 13, then correct B13 (91830 -> 91833)
Swap-Rows 12 13
$ws.Range("B13").Value = 91833

# 4. Full swap of row 14 and row 15
Swap-Rows 14 15

# 5. Full swap of row 18 and row 19
Swap-Rows 18 19

# 6. Three-way rotation: new23 = old24, new24 = old25, new25 = old23
$v23 = Get-RowValues 23
$v24 = Get-RowValues 24
$v25 = Get-RowValues 25
Set-RowValues 23 $v24
Set-RowValues 24 $v25
Set-RowValues 25 $v23

# 7. Full swap of row 29 and row 30
Swap-Rows 29 30
